$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Silan%Gu%NULL%0, Yanfei%Chen%NULL%0, Zhengjie%Wu%NULL%1, Yunbo%Chen%NULL%1, Hainv%Gao%NULL%1, Longxian%Lv%NULL%0, Feifei%Guo%NULL%1, Xuewu%Zhang%NULL%1, Rui%Luo%NULL%1, Chenjie%Huang%NULL%1, Haifeng%Lu%NULL%1, Beiwen%Zheng%NULL%1, Jiaying%Zhang%NULL%1, Ren%Yan%NULL%1, Hua%Zhang%NULL%0, Huiyong%Jiang%NULL%1, Qiaomai%Xu%NULL%1, Jing%Guo%NULL%0, Yiwen%Gong%NULL%1, Lingling%Tang%NULL%1, Lanjuan%Li%ljli@zju.edu.cn%0]'
$ws.Range('I2').Value = ''
$ws.Range('J2').Value = 'Oxford University Press'
$ws.Range('E3').Value = '[Tao%Zuo%NULL%0, Fen%Zhang%NULL%2, Grace C.Y.%Lui%NULL%2, Yun Kit%Yeoh%NULL%1, Amy Y.L.%Li%NULL%1, Hui%Zhan%NULL%2, Yating%Wan%NULL%1, Arthur C.K.%Chung%NULL%1, Chun Pan%Cheung%NULL%1, Nan%Chen%NULL%2, Christopher K.C.%Lai%NULL%1, Zigui%Chen%NULL%1, Eugene Y.K.%Tso%NULL%2, Kitty S.C.%Fung%NULL%1, Veronica%Chan%NULL%1, Lowell%Ling%NULL%1, Gavin%Joynt%NULL%1, David S.C.%Hui%NULL%0, Francis K.L.%Chan%NULL%2, Paul K.S.%Chan%NULL%2, Siew C.%Ng%NULL%2]'
$ws.Range('I3').Value = ''
$ws.Range('J3').Value = 'by the AGA Institute. Published by Elsevier Inc.'
$ws.Range('E4').Value = '[Tao%Zuo%NULL%0, Hui%Zhan%NULL%0, Fen%Zhang%NULL%0, Qin%Liu%NULL%0, Eugene Y.K.%Tso%NULL%0, Grace C.Y.%Lui%NULL%0, Nan%Chen%NULL%0, Amy%Li%NULL%1, Wenqi%Lu%NULL%1, Francis K.L.%Chan%NULL%0, Paul K.S.%Chan%NULL%0, Siew C.%Ng%NULL%0]'
$ws.Range('I4').Value = ''
$ws.Range('J4').Value = 'by the AGA Institute. Published by Elsevier Inc.'
$ws.Range('E5').Value = '[Haocheng%Zhang%NULL%0, Jing-Wen%Ai%NULL%0, Wenjiao%Yang%NULL%1, Xian%Zhou%NULL%0, Fusheng%He%NULL%1, Shumei%Xie%NULL%1, Weiqi%Zeng%NULL%1, Yang%Li%NULL%0, Yiqi%Yu%NULL%0, Xuejing%Gou%NULL%1, Yongjun%Li%NULL%0, Xiaorui%Wang%NULL%1, Hang%Su%NULL%1, Teng%Xu%NULL%0, Wenhong%Zhang%wenhongzhang_hs@126.com%0]'
$ws.Range('I5').Value = ''
$ws.Range('J5').Value = 'Oxford University Press'
$ws.Range('D6').Value = 'id="Par1">We analyzed the bacterial communities of the nasopharynx in 40 SARS-CoV-2 infected and uninfected patients.
 All infected patients had a mild COVID-19 disease.
 We did not find statistically significant differences in either bacterial richness and diversity or composition.
 These findings suggest a nasopharyngeal microbiota at least early resilient to SARS-CoV-2 infection.
'
$ws.Range('E6').Value = '[Flavio%De Maio%NULL%0, Brunella%Posteraro%NULL%1, Francesca Romana%Ponziani%NULL%1, Paola%Cattani%NULL%1, Antonio%Gasbarrini%NULL%0, Maurizio%Sanguinetti%maurizio.sanguinetti@unicatt.it%1]'
$ws.Range('I6').Value = ''
$ws.Range('J6').Value = 'BioMed Central'
$ws.Range('D7').Value = 'From December 2019, an outbreak of unusual pneumonia was reported in Wuhan with many cases linked to Huanan Seafood Market that sells seafood as well as live exotic animals.
 We investigated two patients who developed acute respiratory syndromes after independent contact history with this market.
 The two patients shared common clinical features including fever, cough, and multiple ground-glass opacities in the bilateral lung field with patchy infiltration.
 Here, we highlight the use of a low-input metagenomic next-generation sequencing (mNGS) approach on RNA extracted from bronchoalveolar lavage fluid (BALF).
 It rapidly identified a novel coronavirus (named 2019-nCoV according to World Health Organization announcement) which was the sole pathogens in the sample with very high abundance level (1.5% and 0.62% of total RNA sequenced).
 The entire viral genome is 29,881 nt in length (GenBank MN988668 and MN988669, Sequence Read Archive database Bioproject accession PRJNA601736) and is classified into β-coronavirus genus.
 Phylogenetic analysis indicates that 2019-nCoV is close to coronaviruses (CoVs) circulating in Rhinolophus (Horseshoe bats), such as 98.7% nucleotide identity to partial RdRp gene of bat coronavirus strain BtCoV/4991 (GenBank KP876546, 370 nt sequence of RdRp and lack of other genome sequence) and 87.9% nucleotide identity to bat coronavirus strain bat-SL-CoVZC45 and bat-SL-CoVZXC21. Evolutionary analysis based on ORF1a/1b, S, and N genes also suggests 2019-nCoV is more likely a novel CoV independently introduced from animals to humans.
'
$ws.Range('E7').Value = '[Liangjun%Chen%NULL%0, Weiyong%Liu%NULL%0, Qi%Zhang%NULL%0, Ke%Xu%NULL%2, Ke%Xu%NULL%0, Guangming%Ye%NULL%0, Weichen%Wu%NULL%1, Ziyong%Sun%NULL%0, Fang%Liu%NULL%0, Kailang%Wu%NULL%0, Bo%Zhong%NULL%1, Yi%Mei%NULL%1, Wenxia%Zhang%NULL%1, Yu%Chen%NULL%0, Yirong%Li%NULL%0, Yirong%Li%NULL%0, Mang%Shi%NULL%1, Ke%Lan%NULL%1, Yingle%Liu%NULL%1]'
$ws.Range('I7').Value = ''
$ws.Range('J7').Value = 'Taylor &amp; Francis'
$ws.Range('E8').Value = '[Li-Li%Ren%NULL%0, Ye-Ming%Wang%NULL%0, Zhi-Qiang%Wu%NULL%0, Zi-Chun%Xiang%NULL%0, Li%Guo%NULL%0, Teng%Xu%NULL%0, Yong-Zhong%Jiang%NULL%0, Yan%Xiong%NULL%0, Yong-Jun%Li%NULL%0, Xing-Wang%Li%NULL%0, Hui%Li%NULL%0, Guo-Hui%Fan%NULL%0, Xiao-Ying%Gu%NULL%0, Yan%Xiao%NULL%0, Hong%Gao%NULL%0, Jiu-Yang%Xu%NULL%0, Fan%Yang%NULL%0, Xin-Ming%Wang%NULL%0, Chao%Wu%NULL%0, Lan%Chen%NULL%0, Yi-Wei%Liu%NULL%0, Bo%Liu%NULL%0, Jian%Yang%NULL%0, Xiao-Rui%Wang%NULL%0, Jie%Dong%NULL%0, Li%Li%NULL%0, Chao-Lin%Huang%NULL%0, Jian-Ping%Zhao%NULL%0, Yi%Hu%NULL%0, Zhen-Shun%Cheng%NULL%0, Lin-Lin%Liu%NULL%0, Zhao-Hui%Qian%NULL%0, Chuan%Qin%NULL%0, Qi%Jin%NULL%0, Bin%Cao%NULL%0, Jian-Wei%Wang%NULL%0, Xiu-Yuan%Hao%NULL%0, Xiu-Yuan%Hao%NULL%0, Pei-Fang%Wei%NULL%0]'
$ws.Range('I8').Value = ''
$ws.Range('J8').Value = 'Wolters Kluwer Health'
$ws.Range('E9').Value = '[Jun%Fan%NULL%0, Xiang%Li%NULL%1, Yong%Gao%NULL%0, Junjie%Zhou%NULL%1, Sihua%Wang%NULL%1, Bo%Huang%NULL%1, Junhua%Wu%NULL%1, Qin%Cao%NULL%1, Yajun%Chen%NULL%1, Zhenkao%Wang%NULL%1, Danju%Luo%NULL%1, Ting%Zhou%NULL%0, Ruiting%Li%NULL%0, You%Shang%NULL%0, Xiu%Nie%whunhpath@163.com%0]'
$ws.Range('I9').Value = ''
$ws.Range('J9').Value = 'The British Infection Association. Published by Elsevier Ltd.'
$ws.Range('D10').Value = 'The ongoing global pandemic of COVID-19 disease, which is caused by severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2), mainly infect lung epithelial cells, and spread mainly through respiratory droplets.
 However, recent studies showed potential intestinal infection of SARS-CoV-2, implicated the possibility that the intestinal infection of SARS-CoV-2 may correlate with the dysbiosis of gut microbiota, as well as the severity of COVID-19 symptoms.
 Here, we investigated the alteration of the gut microbiota in COVID-19 patients, as well as analyzed the correlation between the altered microbes and the levels of intestinal inflammatory cytokine IL-18, which was reported to be elevated in the serum of in COVID-19 patients.
 Comparing with healthy controls or seasonal flu patients, the gut microbiota showed significantly reduced diversity, with increased opportunistic pathogens in COVID-19 patients.
 Also, IL-18 level was higher in the fecal samples of COVID-19 patients than in those of either healthy controls or seasonal flu patients.
 Moreover, the IL-18 levels were even higher in the fecal supernatants obtained from COVID-19 patients that tested positive for SARS-CoV-2 RNA than those that tested negative in fecal samples.
 These results indicate that changes in gut microbiota composition might contribute to SARS-CoV-2-induced production of inflammatory cytokines in the intestine and potentially also to the onset of a cytokine storm.
'
$ws.Range('E10').Value = '[Wanyin%Tao%NULL%0, Guorong%Zhang%NULL%1, Xiaofang%Wang%NULL%1, Meng%Guo%NULL%1, Weihong%Zeng%NULL%1, Zhihao%Xu%NULL%1, Dan%Cao%NULL%1, Aijun%Pan%NULL%0, Yucai%Wang%NULL%1, Kaiguang%Zhang%NULL%1, Xiaoling%Ma%NULL%1, Zhengxu%Chen%NULL%1, Tengchuan%Jin%NULL%1, Lianxin%Liu%NULL%0, Jianping%Weng%NULL%1, Shu%Zhu%NULL%1]'
$ws.Range('I10').Value = ''
$ws.Range('J10').Value = 'The Author(s). Published by Elsevier B.V.'
$ws.Range('C11').Value = 'Unknown Title'
$ws.Range('D11').Value = 'Unknown Abstract'
$ws.Range('E11').Value = '[]'
$ws.Range('F11').Value = 'not found'
$ws.Range('G11').Value = 'N/A'
$ws.Range('H11').Value = '1970-01-01'
$ws.Range('I11').Value = ''
$ws.Range('C12').Value = 'Unknown Title'
$ws.Range('D12').Value = 'Unknown Abstract'
$ws.Range('E12').Value = '[]'
$ws.Range('F12').Value = 'not found'
$ws.Range('G12').Value = 'N/A'
$ws.Range('H12').Value = '1970-01-01'
$ws.Range('I12').Value = ''
